$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('https://freelance.habr.com/tasks/594520', 'Внести изменения в код (Django), сделать деплой на AWS', 'django')
    ,@('https://freelance.habr.com/tasks/593886', 'Разработка проектов (ecommerce, аналитика, интеграции) на Django и Fast API', 'python, django, drf, fastapi, sql, postgres')
    ,@('https://freelance.habr.com/tasks/593550', 'Парсер маркетплейсов', 'django, phyton')
    ,@('https://freelance.habr.com/tasks/593374', 'Расширение для wildberries', 'парсинг, django, python')
    ,@('https://freelance.habr.com/tasks/592972', 'Восстановить вебсайт на django', 'django, celery, redis, postgers, nginx, elasticsearch, centos, devops, docker, git')
    ,@('https://freelance.habr.com/tasks/592590', 'Доработка кода back-end django (middle)', 'django, back-end, разработка, python')
    ,@('https://freelance.habr.com/tasks/592062', 'Разработать web приложение для просмотра аниме', 'бэкенд, фронтенд')
    ,@('https://freelance.habr.com/tasks/591756', 'Кэш сервер для API, исправить ошибки и расширить функционал Django', 'django, python')
    ,@('https://freelance.habr.com/tasks/591324', 'Разработка бэкенда туристического портала и PMS', 'python, django, postresql, elasticsearch')
    ,@('https://freelance.habr.com/tasks/590864', 'Доработка SaaS платформы (Python + Django)', 'python, django, эквайринг')
    ,@('https://freelance.habr.com/tasks/590376', 'Поменять монетизацию в проекте', 'django, python')
    ,@('https://freelance.habr.com/tasks/589502', 'Разработка б2б интернет-магазина на Django', 'python, django, fastapi, celery, drf, gunicorn, docker, postgresql, swagger')
    ,@('https://freelance.habr.com/tasks/589226', 'Backend Developer', 'nestsjs, express, ruby, ruby on rails, python, django, react')
    ,@('https://freelance.habr.com/tasks/598363', 'Внести изменения в код (Django), сделать деплой на AWS', 'django')
    ,@('https://freelance.habr.com/tasks/598323', 'Помочь настроить собственные страницы в админке Django в pet-проекте', 'python, django')
    ,@('https://freelance.habr.com/tasks/598145', 'Написать бэкенд на джанго для небольшого pet-проекта', 'бэкенд, rest api, django, backend')
    ,@('https://freelance.habr.com/tasks/598101', 'Помочь установить Postgresql и настроить его с Django на windows', 'postgresql, django')
    ,@('https://freelance.habr.com/tasks/597379', 'Реализовать полнотекстовый поиск на Django и postgresql', 'python, django, postgresql, full text, fat')
    ,@('https://freelance.habr.com/tasks/597325', 'Устранение ошибки CORS Missing Allow Origin в проекте Django+Tornado', 'python, django, tornado, linux')
    ,@('https://freelance.habr.com/tasks/596223', 'Django + Ninja API: Сделать бэкенд для синхронизации сайта с royalty calendar', 'django, ninja api, git')
    ,@('https://freelance.habr.com/tasks/594821', 'Backend dev для Solana MPC Wallet', 'web3, sss, node.js, phyton, rust')
    ,@('https://freelance.habr.com/tasks/593362', 'Рефакторинг кода backend на Django + деплой на сервер', 'django, python, docker')
    ,@('https://freelance.habr.com/tasks/598325', 'Разработать на DRF отображение древовидных данных', 'django, drf, postgresql, restapi')
    ,@('https://freelance.habr.com/tasks/597721', 'Доработать сервис (только индивидуальные разработчики)', 'python, django')
    ,@('https://freelance.habr.com/tasks/597487', 'Доработка сервиса за прошлой командой', 'python, django')
    ,@('https://freelance.habr.com/tasks/597281', 'Разработка проектов (ecommerce) на Django', 'django, python, ооп')
    ,@('https://freelance.habr.com/tasks/596761', 'Нужно проверить и доработать API в проекте django drf', 'django, django rest framework')
    ,@('https://freelance.habr.com/tasks/595465', 'Доработка и тех. поддержка проекта на Django', 'python, django')
    ,@('https://freelance.habr.com/tasks/595367', 'Разработка на Django', 'django, python')
    ,@('https://freelance.habr.com/tasks/593198', 'Заменить Sendpulse на аналог sendsay, проект на python (django)', 'python, django')
    ,@('https://freelance.habr.com/tasks/555303', 'Доработать существующий код, пофиксить баги, написать новый', 'django ninja, django')
    ,@('https://freelance.habr.com/tasks/446435', 'Разработать бэкэнд для сайта выборов преподавательского состава', 'postgresql, redis, django, nginx, docker, linux, python')
    ,@('https://freelance.habr.com/tasks/598517', 'Разработать telegram бота для управления и продажи VPN', 'python, telegram, vless, vpn, bot, marzban, xray')
    ,@('https://freelance.habr.com/tasks/596965', 'Разработать телеграм-бота с внедренным ИИ', 'telegram, bot, ai')
    ,@('https://freelance.habr.com/tasks/595407', 'Разработка магазин TG-bot', 'python, tg, опыт, ответственность')
    ,@('https://freelance.habr.com/tasks/595303', 'Собрать телеграм-бота из нескольких вопросов для ниши стройка и связать с амо на nocode', 'amocrm, bot, salesbot, salebot')
    ,@('https://freelance.habr.com/tasks/594887', 'Разработать TG-бота для анализа и мониторинга кошельков в Ethereum', 'telegram, blockchain, ethereum, telegram bot')
    ,@('https://freelance.habr.com/tasks/593832', 'Консультация TELEGRAM API', 'js, node, telegram api, mtproto​, bot')
    ,@('https://freelance.habr.com/tasks/593672', 'Разработка TON-площадки в Telegram Mini App', 'блокчейн ton, telegram mini app, web3')
    ,@('https://freelance.habr.com/tasks/593446', 'Разработать чат-бота для веб-сайта', 'python, bot, chatgpt')
    ,@('https://freelance.habr.com/tasks/593388', 'Telegram Bot для проверки постбеков', 'телеграм, бот, постбек')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
}

Write-Host "Done writing rows"